$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p084v_a4</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p084v_4</id>", 2)

$d.Content.Find.Execute("<id>p085r_a1</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p085r_1</id>", 2)

$d.Content.Find.Execute("<id>p085r_a2</id>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p085r_2</id>", 2)
